# Revert "Drop in all data files from 3.0 RMI script"
# Adds a new "Texas Notes" worksheet (between "About" and "PPEIdtIL") that
# contains review notes/questions about the 10% vs 2% labeling assumption.

$wb = $excel.ActiveWorkbook

# --- Insert the new "Texas Notes" sheet between "About" and "PPEIdtIL" ---
$ppeSheet = $wb.Worksheets.Item("PPEIdtIL")
$notes = $wb.Worksheets.Add($ppeSheet)
$notes.Name = "Texas Notes"

# --- Populate "Texas Notes" content ---
$notes.Range("A1").Value = "The study says:"

$notes.Range("A3").Value = "Table 2. Costs and benefits of improved appliance labeling"
$notes.Range("A4").Value = "Average % energy savings"

$notes.Range("A5").Value = 0.1
$notes.Range("A5").NumberFormat = "0%"

# Row 6-10: merged note block, left/top aligned, wrapped
$notes.Range("A6:M10").Merge()
$notes.Range("A6").Value = "Assumes consumer selects products that average 10% energy savings. Difference in energy use from least to most efficient product varies from 10% to 50% depending on product category according to FTC published ranges. We assume most consumers motivated by the label select a product that falls 1-2 categories more efficient (e.g., shift from 1-star to 2- or 3-star product, or from 3-star to 4- or 5-star product, etc.). Source: ACEEE estimate based on labeling research."
$notes.Range("A6:M10").HorizontalAlignment = -4131
$notes.Range("A6:M10").VerticalAlignment = -4160
$notes.Range("A6:M10").WrapText = $true

# Row 11: header row, left/top aligned (no wrap); fillers additionally wrapped
$notes.Range("B11:M11").HorizontalAlignment = -4131
$notes.Range("B11:M11").VerticalAlignment = -4160
$notes.Range("B11:M11").WrapText = $true
$notes.Range("A11").Value = "% of market affected by policy"
$notes.Range("A11").HorizontalAlignment = -4131
$notes.Range("A11").VerticalAlignment = -4160
$notes.Range("A11").WrapText = $false

$notes.Range("A12").Value = 0.2
$notes.Range("A12").NumberFormat = "0%"

# Row 13-14: merged note block, left/top aligned, wrapped
$notes.Range("A13:M14").Merge()
$notes.Range("A13").Value = "Assumes 20% of consumers are motivated by the improved label and use it to select more efficient products. Source: ACEEE estimate based on prior labeling research"
$notes.Range("A13:M14").HorizontalAlignment = -4131
$notes.Range("A13:M14").VerticalAlignment = -4160
$notes.Range("A13:M14").WrapText = $true

# Row 15: blank spacer row, top aligned + wrapped (general horizontal)
$notes.Range("A15:M15").HorizontalAlignment = 1
$notes.Range("A15:M15").VerticalAlignment = -4160
$notes.Range("A15:M15").WrapText = $true

# Rows 16-17: top aligned only (general horizontal, no wrap)
$notes.Range("A16").Value = "Question for EI"
$notes.Range("A16:M16").HorizontalAlignment = 1
$notes.Range("A16:M16").VerticalAlignment = -4160
$notes.Range("A16:M16").WrapText = $false

$notes.Range("A17").Value = "does this mean that the numbers in this spreadsheet should be 2% (10% * 20%)?"
$notes.Range("A17:M17").HorizontalAlignment = 1
$notes.Range("A17:M17").VerticalAlignment = -4160
$notes.Range("A17:M17").WrapText = $false

# Rows 18-23: default/general alignment (just flagged, no actual overrides)
$notes.Range("A18:M23").HorizontalAlignment = 1
$notes.Range("A18:M23").VerticalAlignment = -4107
$notes.Range("A18:M23").WrapText = $false
$notes.Range("A18").Value = "right now they are 10%."
$notes.Range("A20").Value = "*Anwer from EI: yes, this was an error."

# --- View/selection tweaks to mirror the target workbook state ---
# (Re-fetch sheet references by name: the worksheet-collection insert above
# can leave previously-held references pointing at the wrong sheet object.)
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Range("D28").Select()

$ppeSheet = $wb.Worksheets.Item("PPEIdtIL")
$ppeSheet.Activate()
$ppeSheet.Range("D10").Select()

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
$aboutSheet.Range("C24").Select()
